$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (existing data ends at row 89)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

# Row 90: G1 / Test1
$ws.Cells.Item($newRow1, 1).Value = "G1"
$ws.Cells.Item($newRow1, 2).Value = "Test1"
$ws.Cells.Item($newRow1, 3).Value = 45905
$ws.Cells.Item($newRow1, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow1, 4).Value = 0.6518999194026193
$ws.Cells.Item($newRow1, 5).Value = 0
$ws.Cells.Item($newRow1, 6).Value = -0.01

# Row 91: G2 / sedrftgyhuioygtfrd
$ws.Cells.Item($newRow2, 1).Value = "G2"
$ws.Cells.Item($newRow2, 2).Value = "sedrftgyhuioygtfrd"
$ws.Cells.Item($newRow2, 3).Value = 45905
$ws.Cells.Item($newRow2, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow2, 4).Value = 0.6518999194026193
$ws.Cells.Item($newRow2, 5).Value = 0
$ws.Cells.Item($newRow2, 6).Value = -0.01
